# Update the "取得日時" (acquired datetime) column for rows 2-17 on the
# "ランサーズ" sheet from 2025-10-14 18:25:30 to 2025-10-14 18:34:10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-14 18:34:10"

for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
